$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Known blank inline-string cell used as a template for empty placeholder cells
$blank = $ws.Range("I2")

# ---- Row 3 ----
$ws.Cells.Item(3, 1).Value = 112394252
$ws.Cells.Item(3, 2).Value = 89557
$ws.Cells.Item(3, 3).Value = "Ovaliderad"
$ws.Cells.Item(3, 4).Value = "NT"
$ws.Cells.Item(3, 5).Value = 5432
$ws.Cells.Item(3, 6).Value = "Granticka"
$ws.Cells.Item(3, 7).Value = "Porodaedalea chrysoloma"
$ws.Cells.Item(3, 8).Value = "(Fr.) Fiasson & Niemelä"
$blank.Copy($ws.Cells.Item(3, 9))
$ws.Cells.Item(3, 10).Value = "fruktkroppar"
$blank.Copy($ws.Cells.Item(3, 11))
$blank.Copy($ws.Cells.Item(3, 14))
$ws.Cells.Item(3, 16).Value = "Paurankivägen, T lm"
$ws.Cells.Item(3, 17).Value = 749928
$ws.Cells.Item(3, 18).Value = 7535991
$ws.Cells.Item(3, 19).Value = 10
$ws.Cells.Item(3, 20).Value = "Norrbotten"
$ws.Cells.Item(3, 21).Value = "Kiruna"
$ws.Cells.Item(3, 22).Value = "Torne lappmark"
$ws.Cells.Item(3, 23).Value = "Jukkasjärvi"
$ws.Cells.Item(3, 25).NumberFormat = "@"
$ws.Cells.Item(3, 25).Value = "2023-09-28"
$ws.Cells.Item(3, 25).ClearFormats()
$ws.Cells.Item(3, 27).NumberFormat = "@"
$ws.Cells.Item(3, 27).Value = "2023-09-28"
$ws.Cells.Item(3, 27).ClearFormats()
$ws.Cells.Item(3, 30).Value = $false
$ws.Cells.Item(3, 31).Value = $false
$blank.Copy($ws.Cells.Item(3, 32))
$ws.Cells.Item(3, 33).Value = $false
$blank.Copy($ws.Cells.Item(3, 46))
$ws.Cells.Item(3, 49).Value = "Anne Järvinen"
$ws.Cells.Item(3, 50).Value = "Anne Järvinen"
$blank.Copy($ws.Cells.Item(3, 51))

# ---- Row 4 ----
$ws.Cells.Item(4, 1).Value = 112394378
$ws.Cells.Item(4, 2).Value = 89557
$ws.Cells.Item(4, 3).Value = "Ovaliderad"
$ws.Cells.Item(4, 4).Value = "NT"
$ws.Cells.Item(4, 5).Value = 5432
$ws.Cells.Item(4, 6).Value = "Granticka"
$ws.Cells.Item(4, 7).Value = "Porodaedalea chrysoloma"
$ws.Cells.Item(4, 8).Value = "(Fr.) Fiasson & Niemelä"
$blank.Copy($ws.Cells.Item(4, 9))
$ws.Cells.Item(4, 10).Value = "fruktkroppar"
$blank.Copy($ws.Cells.Item(4, 11))
$blank.Copy($ws.Cells.Item(4, 14))
$ws.Cells.Item(4, 16).Value = "Paurankivägen, T lm"
$ws.Cells.Item(4, 17).Value = 749921
$ws.Cells.Item(4, 18).Value = 7535985
$ws.Cells.Item(4, 19).Value = 10
$ws.Cells.Item(4, 20).Value = "Norrbotten"
$ws.Cells.Item(4, 21).Value = "Kiruna"
$ws.Cells.Item(4, 22).Value = "Torne lappmark"
$ws.Cells.Item(4, 23).Value = "Jukkasjärvi"
$ws.Cells.Item(4, 25).NumberFormat = "@"
$ws.Cells.Item(4, 25).Value = "2023-09-28"
$ws.Cells.Item(4, 25).ClearFormats()
$ws.Cells.Item(4, 27).NumberFormat = "@"
$ws.Cells.Item(4, 27).Value = "2023-09-28"
$ws.Cells.Item(4, 27).ClearFormats()
$ws.Cells.Item(4, 30).Value = $false
$ws.Cells.Item(4, 31).Value = $false
$blank.Copy($ws.Cells.Item(4, 32))
$ws.Cells.Item(4, 33).Value = $false
$blank.Copy($ws.Cells.Item(4, 46))
$ws.Cells.Item(4, 49).Value = "Anne Järvinen"
$ws.Cells.Item(4, 50).Value = "Anne Järvinen"
$blank.Copy($ws.Cells.Item(4, 51))

# ---- Row 5 ----
$ws.Cells.Item(5, 1).Value = 112394151
$ws.Cells.Item(5, 2).Value = 89820
$ws.Cells.Item(5, 3).Value = "Ovaliderad"
$ws.Cells.Item(5, 4).Value = "NT"
$ws.Cells.Item(5, 5).Value = 658
$ws.Cells.Item(5, 6).Value = "Rosenticka"
$ws.Cells.Item(5, 7).Value = "Rhodofomes roseus"
$ws.Cells.Item(5, 8).Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$blank.Copy($ws.Cells.Item(5, 9))
$ws.Cells.Item(5, 10).Value = "fruktkroppar"
$blank.Copy($ws.Cells.Item(5, 11))
$blank.Copy($ws.Cells.Item(5, 14))
$ws.Cells.Item(5, 16).Value = "Paurankivägen, T lm"
$ws.Cells.Item(5, 17).Value = 749936
$ws.Cells.Item(5, 18).Value = 7535990
$ws.Cells.Item(5, 19).Value = 10
$ws.Cells.Item(5, 20).Value = "Norrbotten"
$ws.Cells.Item(5, 21).Value = "Kiruna"
$ws.Cells.Item(5, 22).Value = "Torne lappmark"
$ws.Cells.Item(5, 23).Value = "Jukkasjärvi"
$ws.Cells.Item(5, 25).NumberFormat = "@"
$ws.Cells.Item(5, 25).Value = "2023-09-28"
$ws.Cells.Item(5, 25).ClearFormats()
$ws.Cells.Item(5, 27).NumberFormat = "@"
$ws.Cells.Item(5, 27).Value = "2023-09-28"
$ws.Cells.Item(5, 27).ClearFormats()
$ws.Cells.Item(5, 30).Value = $false
$ws.Cells.Item(5, 31).Value = $false
$blank.Copy($ws.Cells.Item(5, 32))
$ws.Cells.Item(5, 33).Value = $false
$blank.Copy($ws.Cells.Item(5, 46))
$ws.Cells.Item(5, 49).Value = "Anne Järvinen"
$ws.Cells.Item(5, 50).Value = "Anne Järvinen"
$blank.Copy($ws.Cells.Item(5, 51))

# ---- Row 6 ----
$ws.Cells.Item(6, 1).Value = 112457286
$ws.Cells.Item(6, 2).Value = 89820
$ws.Cells.Item(6, 3).Value = "Ovaliderad"
$ws.Cells.Item(6, 4).Value = "NT"
$ws.Cells.Item(6, 5).Value = 658
$ws.Cells.Item(6, 6).Value = "Rosenticka"
$ws.Cells.Item(6, 7).Value = "Rhodofomes roseus"
$ws.Cells.Item(6, 8).Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$blank.Copy($ws.Cells.Item(6, 9))
$blank.Copy($ws.Cells.Item(6, 10))
$blank.Copy($ws.Cells.Item(6, 11))
$blank.Copy($ws.Cells.Item(6, 14))
$ws.Cells.Item(6, 16).Value = "Paurankivaravägen, T lm"
$ws.Cells.Item(6, 17).Value = 749868
$ws.Cells.Item(6, 18).Value = 7536000
$ws.Cells.Item(6, 19).Value = 10
$ws.Cells.Item(6, 20).Value = "Norrbotten"
$ws.Cells.Item(6, 21).Value = "Kiruna"
$ws.Cells.Item(6, 22).Value = "Torne lappmark"
$ws.Cells.Item(6, 23).Value = "Jukkasjärvi"
$ws.Cells.Item(6, 25).NumberFormat = "@"
$ws.Cells.Item(6, 25).Value = "2023-10-01"
$ws.Cells.Item(6, 25).ClearFormats()
$ws.Cells.Item(6, 27).NumberFormat = "@"
$ws.Cells.Item(6, 27).Value = "2023-10-01"
$ws.Cells.Item(6, 27).ClearFormats()
$ws.Cells.Item(6, 30).Value = $false
$ws.Cells.Item(6, 31).Value = $false
$blank.Copy($ws.Cells.Item(6, 32))
$ws.Cells.Item(6, 33).Value = $false
$blank.Copy($ws.Cells.Item(6, 46))
$ws.Cells.Item(6, 49).Value = "Stefan Andersson"
$ws.Cells.Item(6, 50).Value = "Stefan Andersson, per-erik mukka, Christina Boyd"
$blank.Copy($ws.Cells.Item(6, 51))

# ---- Row 7 ----
$ws.Cells.Item(7, 1).Value = 112438309
$ws.Cells.Item(7, 2).Value = 90434
$ws.Cells.Item(7, 3).Value = "Ovaliderad"
$ws.Cells.Item(7, 4).Value = "NT"
$ws.Cells.Item(7, 5).Value = 4745
$ws.Cells.Item(7, 6).Value = "Tallriska"
$ws.Cells.Item(7, 7).Value = "Lactarius musteus"
$ws.Cells.Item(7, 8).Value = "Fr."
$blank.Copy($ws.Cells.Item(7, 9))
$blank.Copy($ws.Cells.Item(7, 10))
$blank.Copy($ws.Cells.Item(7, 11))
$blank.Copy($ws.Cells.Item(7, 14))
$ws.Cells.Item(7, 16).Value = "Paurankivaravägen, T lm"
$ws.Cells.Item(7, 17).Value = 749892
$ws.Cells.Item(7, 18).Value = 7535980
$ws.Cells.Item(7, 19).Value = 25
$ws.Cells.Item(7, 20).Value = "Norrbotten"
$ws.Cells.Item(7, 21).Value = "Kiruna"
$ws.Cells.Item(7, 22).Value = "Torne lappmark"
$ws.Cells.Item(7, 23).Value = "Jukkasjärvi"
$ws.Cells.Item(7, 25).NumberFormat = "@"
$ws.Cells.Item(7, 25).Value = "2023-10-01"
$ws.Cells.Item(7, 25).ClearFormats()
$ws.Cells.Item(7, 26).Value = "11:54"
$ws.Cells.Item(7, 27).NumberFormat = "@"
$ws.Cells.Item(7, 27).Value = "2023-10-01"
$ws.Cells.Item(7, 27).ClearFormats()
$ws.Cells.Item(7, 28).Value = "11:54"
$ws.Cells.Item(7, 30).Value = $false
$ws.Cells.Item(7, 31).Value = $false
$blank.Copy($ws.Cells.Item(7, 32))
$ws.Cells.Item(7, 33).Value = $false
$blank.Copy($ws.Cells.Item(7, 46))
$ws.Cells.Item(7, 49).Value = "Stefan Andersson"
$ws.Cells.Item(7, 50).Value = "Stefan Andersson, Christina Boyd, per-erik mukka"
$blank.Copy($ws.Cells.Item(7, 51))

$excel.CutCopyMode = $false
